# Gestion projet S2.xlsx - apply the edits described by the commit
#   "update color black to grey update GUI F11 key for fullscreen update Gestion projet"
# The substantive part of this commit updates the weekly tracking table on
# "Feuil1": a few hour/cost entries are corrected for weeks 15-17 (rows 17-19)
# and a missing running-total formula is restored for week 16 (row 18). The
# dependent formula cells (J/N/O/V/W/X columns, and the W24/G25/H25/I25
# summaries) recalculate automatically; only the literal inputs are set here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 7: Prestation expert prévu (S) now explicitly recorded as 0
$ws.Range("S7").Value = 0

# Row 11: Location materiel (U) effective cost recorded
$ws.Range("U11").Value = 100

# Row 17 (week 15): Nicolas/Simon organisation effectif corrected
$ws.Range("K17").Value = 4
$ws.Range("M17").Value = 4

# Row 18 (week 16): technique + organisation effectif hours filled in, plus
# the prevu/effectif toggle between "Location locaux" (S) and "Prestation
# expert" (T)
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1
$ws.Range("K18").Value = 9
$ws.Range("L18").Value = 5
$ws.Range("M18").Value = 11
$ws.Range("Q18").Value = 15
$ws.Range("S18").Value = 100
$ws.Range("T18").ClearContents()

# Row 18's running-total formula (X18) was missing; restore it so it keeps
# pace with the other weeks (X11..X17, X19..)
$ws.Range("X18").Formula = "=SUM(`$W`$3:W18)+SUM(V19:`$V`$23)"

# Row 19 (week 17): organisation effectif filled in, prevu Location locaux
# (T) recorded
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 1
$ws.Range("T19").Value = 100

# Keep the saved cursor position in sync with where the edit left off
$ws.Range("X22").Select()
